$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (date 25.12.2024 -> 30.12.2024) ---
# Order matters for shared-string table layout: write A10 before A1
# so the new strings land at sst indices 42/43 matching the source workbook.
$ws.Range("A10").Value = "DAILY STOCK                         (30/12/2024) "
$ws.Range("A1").Value = "Mangrove Communication   30.12.2024"

# --- Sale/Due table (rows 3-6) updates ---
$ws.Range("C3").Value = 23361
$ws.Range("D3").Value = 27231

$ws.Range("C4").Value = 20776
$ws.Range("D4").Value = 5550
$ws.Range("F4").Value = $null

$ws.Range("C5").Value = 17687
$ws.Range("D5").Value = 19989

$ws.Range("C6").Value = 23223
$ws.Range("D6").Value = 1110

# --- Stock table updates ---
$ws.Range("C13").Value = 157384

$ws.Range("C14").Value = 264450
$ws.Range("D14").Value = 85087
$ws.Range("E14").Value = 207792

$ws.Range("C20").Value = 4150
$ws.Range("D20").Value = 2760

$ws.Range("C21").Value = 220
$ws.Range("D21").Value = 100

$ws.Range("C24").Value = 12

$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 5

$ws.Range("C27").Value = 75
$ws.Range("D27").Value = 1

# --- Selection update ---
$ws.Range("D28").Select()
